$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bishal's image URL (row 7, column B) with the new meme image
$ws.Range("B7").Value = "meme_images/bishal_meme.png"

# Update the active selection to match the saved workbook view
$ws.Range("I12").Select()
